$d = $word.ActiveDocument
$p3 = $d.Paragraphs.Item(3)
$p4 = $d.Paragraphs.Item(4)
$delRange = $d.Range($p3.Range.Start, $p4.Range.Start)
$delRange.Delete()

$p3 = $d.Paragraphs.Item(3)
$bmRange = $d.Range($p3.Range.Start, $p3.Range.Start)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# Set language on trailing empty paragraphs (4,5,6)
for ($idx = 4; $idx -le 6; $idx++) {
    $p = $d.Paragraphs.Item($idx)
    $p.Range.LanguageID = "fr-FR"
}

Write-Output "done"
